$d = $word.ActiveDocument

# Helper: insert a new plain paragraph (no style / no list) after the
# current last paragraph, optionally setting its text, and return it.
function Add-PlainParagraph($text) {
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $newP = $d.Paragraphs.Last
    # Strip any inherited numbering / list style from the split so the
    # new paragraph starts out completely bare, like the source diff.
    $newP.Range.ListFormat.RemoveNumbers()
    $newP.Style = "Normal"
    if ($text) {
        $newP.Range.Text = $text
    }
    return $newP
}

# Blank separator paragraph.
Add-PlainParagraph $null | Out-Null

# Date heading.
Add-PlainParagraph "3/8/22" | Out-Null

# Section heading.
Add-PlainParagraph "Main Features:" | Out-Null

# First bullet item - starts a brand-new bulleted list (new numId).
$item1 = Add-PlainParagraph "Added a service component that handles http services."
$item1.Style = "List Paragraph"
$bulletTemplate = $word.ListGalleries.Item(1).ListTemplates.Item(1)
$item1.Range.ListFormat.ApplyListTemplate($bulletTemplate)

# Second bullet item - continues the same list/numId as the first item.
$item2 = Add-PlainParagraph "Put mock data into firebase database and linked the data into the necessary spots within the web app."
$item2.Style = "List Paragraph"
$item2.Range.ListFormat.List = $item1.Range.ListFormat.List
$item2.Range.ListFormat.ListLevelNumber = $item1.Range.ListFormat.ListLevelNumber

# Trailing blank paragraph, indented to match the end of the list block.
$trailing = Add-PlainParagraph $null
$trailing.LeftIndent = 18
